$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 161, shifting the old rows 161-172
# down to 162-173.
$ws.Rows.Item(161).Insert()

# Populate the new row 161 with a new weekly price observation. Columns
# A,B,C,E,F,G,H,I,J,K,L,Q,T carry the same constant values as the other
# rows in this Mango/Vega Monumental Concepcion block.
$ws.Cells.Item(161, 1).Value = 11
$ws.Cells.Item(161, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(161, 3).Value = "Bíobío"
$ws.Cells.Item(161, 4).Value = 45106
$ws.Cells.Item(161, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(161, 5).Value = 8
$ws.Cells.Item(161, 6).Value = "Fruta"
$ws.Cells.Item(161, 7).Value = 100108
$ws.Cells.Item(161, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(161, 9).Value = 100108002
$ws.Cells.Item(161, 10).Value = "Mango"
$ws.Cells.Item(161, 11).Value = "Sin especificar"
$ws.Cells.Item(161, 12).Value = "Primera"
$ws.Cells.Item(161, 13).Value = 110
$ws.Cells.Item(161, 14).Value = 9000
$ws.Cells.Item(161, 15).Value = 10000
$ws.Cells.Item(161, 16).Value = 9455
$ws.Cells.Item(161, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(161, 18).Value = "Perú"
$ws.Cells.Item(161, 19).Value = 2364
$ws.Cells.Item(161, 20).Value = 4
